$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price/Volume columns so values like "1.00" or
# "9.00" are not auto-coerced to numbers by Excel's type inference, matching
# the original inline-string ("text") cell contents.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

# --- Update Price (D) and Volume(1h) (E) columns for changed rows ---
$ws.Range("D2").Value = "69.092.29"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "3.774.30"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "624.87"
$ws.Range("E5").Value = "  +4.33%  "
$ws.Range("D6").Value = "166.58"
$ws.Range("E6").Value = "  +2.23%  "
$ws.Range("D7").Value = "3.772.94"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  +1.67%  "
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("D11").Value = "0.458"
$ws.Range("E11").Value = "  +3.19%  "
$ws.Range("D12").Value = "6.71"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("D14").Value = "35.69"
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").Value = "4.411.95"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "3.865.02"
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("D17").Value = "69.094.74"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").Value = "17.66"
$ws.Range("E18").Value = "  -2.70%  "
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").Value = "468.43"
$ws.Range("E21").Value = "  +2.67%  "
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("E23").Value = "  +2.54%  "
$ws.Range("E24").Value = "  +4.73%  "
$ws.Range("D25").Value = "83.13"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("D26").Value = "12.05"
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("E27").Value = "  +3.89%  "
$ws.Range("D28").Value = "10.03"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "3.922.85"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").Value = "2.68"
$ws.Range("D32").Value = "2.25"
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D34").Value = "28.77"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D38").Value = "0.164"
$ws.Range("E38").Value = "  +14.10%  "
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("D40").Value = "3.44"
$ws.Range("E40").Value = "  +8.86%  "
$ws.Range("D41").Value = "5.82"
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("D42").Value = "0.966"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D45").Value = "0.298"
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("D46").Value = "43.23"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").Value = "152.49"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("E48").Value = "  +4.48%  "
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("E51").Value = "  -0.05%  "

# --- Rows 36 and 37 swap: Aptos and RenzoRestakedETH switch positions with updated values ---
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "9.00"
$ws.Range("E36").Value = "  +0.90%  "

$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.724.67"
$ws.Range("E37").Value = "  +0.05%  "

# Restore default (unstyled) cell style now that text-typing has been locked in,
# so no stray NumberFormat is left behind on the data cells.
$dataRange.Style = "Normal"

